$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: years 2007-2012 with updated values, replacing old 2003-2012 data
$years = @(2007, 2008, 2009, 2010, 2011, 2012)
$vals  = @(629, 646, 499, 752, 548, 675)

for ($i = 0; $i -lt $years.Length; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $years[$i]
    $ws.Cells.Item($r, 2).Value = $vals[$i]
}

# Remove the now-unused trailing rows (previously rows 7-10)
$ws.Range("A7:B10").ClearContents()

# Update the chart: style and series source ranges now only cover 6 rows
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$chart.ChartStyle = 8

$ser = $chart.SeriesCollection(1)
$ser.Formula = "=SERIES(,'Sheet1'!`$A`$1:`$A`$6,'Sheet1'!`$B`$1:`$B`$6,1)"
